$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-11-15 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-11-16 Sunday", 2) | Out-Null
$d.Content.Find.Execute("735÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "717÷9=", 2) | Out-Null
$d.Content.Find.Execute("991÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "632÷7=", 2) | Out-Null
$d.Content.Find.Execute("467÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "169÷4=", 2) | Out-Null
$d.Content.Find.Execute("494÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "852÷9=", 2) | Out-Null
$d.Content.Find.Execute("115÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "473÷8=", 2) | Out-Null
$d.Content.Find.Execute("648÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "531÷4=", 2) | Out-Null
$d.Content.Find.Execute("191÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "284÷2=", 2) | Out-Null
$d.Content.Find.Execute("727÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "431÷9=", 2) | Out-Null
$d.Content.Find.Execute("959÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "946÷8=", 2) | Out-Null
$d.Content.Find.Execute("238÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "804÷9=", 2) | Out-Null
$d.Content.Find.Execute("116÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "976÷6=", 2) | Out-Null
$d.Content.Find.Execute("882÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "951÷8=", 2) | Out-Null
$d.Content.Find.Execute("532÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "169÷2=", 2) | Out-Null
$d.Content.Find.Execute("493÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "904÷8=", 2) | Out-Null
$d.Content.Find.Execute("507÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "694÷2=", 2) | Out-Null
$d.Content.Find.Execute("531÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "863÷5=", 2) | Out-Null
$d.Content.Find.Execute("272÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "724÷2=", 2) | Out-Null
$d.Content.Find.Execute("324÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "116÷8=", 2) | Out-Null
$d.Content.Find.Execute("844÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "861÷3=", 2) | Out-Null
$d.Content.Find.Execute("753÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "217÷3=", 2) | Out-Null
$d.Content.Find.Execute("763÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "940÷6=", 2) | Out-Null
$d.Content.Find.Execute("703÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "202÷3=", 2) | Out-Null
$d.Content.Find.Execute("890÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "834÷2=", 2) | Out-Null
$d.Content.Find.Execute("568÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "507÷9=", 2) | Out-Null
$d.Content.Find.Execute("288÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "379÷8=", 2) | Out-Null
